$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '48.103.23'
$ws.Range('E2').Value = '  +1.72%  '

$ws.Range('D3').Value = '2.506.03'
$ws.Range('E3').Value = '  +0.62%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '321.27'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.11%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '108.39'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.09%  '

$ws.Range('E7').Value = '  +0.92%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.05%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.542'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.15%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '39.80'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +1.78%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '20.16'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +9.95%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0817'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.85%  '

$ws.Range('E13').Value = '  +0.55%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.19'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.97%  '

$ws.Range('D15').Value = '2.898.89'
$ws.Range('E15').Value = '  +0.65%  '

$ws.Range('D16').Value = '2.502.94'
$ws.Range('E16').Value = '  +0.39%  '

$ws.Range('E17').Value = '  +0.13%  '

$ws.Range('D18').Value = '47.951.96'
$ws.Range('E18').Value = '  +1.62%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.15'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.07%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.64'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.40%  '

$ws.Range('D21').Value = '0.0₃0943'
$ws.Range('E21').Value = '  +0.79%  '

$ws.Range('E22').Value = '  +0.88%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '72.07'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +2.48%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '276.11'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +12.59%  '

$ws.Range('E25').Value = '  +0.13%  '

$ws.Range('E26').Value = '  +0.03%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '25.89'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.77%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.26'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.77%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '10.03'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.64%  '

$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '35.44'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +2.47%  '

$ws.Range('B31').Value = 'Kaspa'
$ws.Range('C31').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.138'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.05%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '49.46'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.70%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '19.36'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -4.42%  '

$ws.Range('E34').Value = '  +0.25%  '

$ws.Range('E35').Value = '  -0.04%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0785'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.25%  '

$ws.Range('E37').Value = '  -0.10%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.59'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -3.43%  '

$ws.Range('E39').Value = '  +1.23%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '122.91'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +3.99%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.111'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.28%  '

$ws.Range('E42').Value = '  -0.99%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '21.74'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -4.99%  '

$ws.Range('E44').Value = '  +3.45%  '

$ws.Range('D45').Value = '2.001.44'
$ws.Range('E45').Value = '  +0.33%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.13'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +3.23%  '

$ws.Range('E47').Value = '  +4.13%  '

$ws.Range('E48').Value = '  -0.44%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '9.02'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -1.20%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '5.18'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.73%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '79.79'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +2.61%  '

